# Generate Report for Handoff
# Updates the localization-status workbook to reflect that b.md has been
# handed off for zh-cn and de-de (new handoff file + datetime), and the
# Overview sheet's status for b.md changes from
# "Handed back: in sync with en-US" to "Ready for handoff".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Overview sheet: update Status for the "b.md" row (row 3) in both
#    the zh-cn (B) and de-de (C) columns.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

# ---------------------------------------------------------------------
# 2. zh-cn sheet: update the "b.md" row (row 3):
#      - Status -> "Ready for handoff"
#      - Latest Handoff File -> new file name (with matching hyperlink text)
#      - Latest Handoff Datetime -> new timestamp
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("B3").Value = "Ready for handoff"
$wsZh.Range("D3").Value = "2016-03-10 04:54:00"

# Rebuild the sheet's hyperlinks so the C3 hyperlink shows the new file
# name while every link keeps pointing at its original target address.
$zhC2Address = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c55a5d5823f8f8dbc47992bd77ffa22f1f728db6/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhF2Address = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/3ca2cd945e27fede4740670e996567560dca847a/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"

$wsZh.Range("C3").Hyperlinks.Delete()

$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/bd6510893bbfe8f56a0bd73bd6287d6a6ed857e3/e2e/a.md", "", "", "a.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("C2"), $zhC2Address, "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/33d1ea1a29d5d295b77e359618e14a9d05af4399/e2e/a.md", "", "", "a.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), $zhF2Address, "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/bd6510893bbfe8f56a0bd73bd6287d6a6ed857e3/e2e/b.md", "", "", "b.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("C3"), $zhC2Address, "", "", "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/33d1ea1a29d5d295b77e359618e14a9d05af4399/e2e/a.md", "", "", "a.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), $zhF2Address, "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/bd6510893bbfe8f56a0bd73bd6287d6a6ed857e3/.localization-config", "", "", ".localization-config") | Out-Null

# Keep the re-created hyperlink cell using the sheet's existing hyperlink look.
$wsZh.Range("C3").Style = "HyperLink"

# ---------------------------------------------------------------------
# 3. de-de sheet: same kind of update for the "b.md" row (row 3).
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("B3").Value = "Ready for handoff"
$wsDe.Range("D3").Value = "2016-03-10 04:54:08"

$deC2Address = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ff812c58013cb17b643962b9fbe98649f1f7bc7d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$deF2Address = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/28b9f4cd7fd9c326da605e6e56086e5d62facf8f/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

$wsDe.Range("C3").Hyperlinks.Delete()

$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/bd6510893bbfe8f56a0bd73bd6287d6a6ed857e3/e2e/a.md", "", "", "a.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("C2"), $deC2Address, "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/b1e5b35d9e5a578fa16bc73772d885e1c2a8ed20/e2e/a.md", "", "", "a.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), $deF2Address, "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/bd6510893bbfe8f56a0bd73bd6287d6a6ed857e3/e2e/b.md", "", "", "b.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("C3"), $deC2Address, "", "", "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/b1e5b35d9e5a578fa16bc73772d885e1c2a8ed20/e2e/a.md", "", "", "a.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), $deF2Address, "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/bd6510893bbfe8f56a0bd73bd6287d6a6ed857e3/.localization-config", "", "", ".localization-config") | Out-Null

# Keep the re-created hyperlink cell using the sheet's existing hyperlink look.
$wsDe.Range("C3").Style = "HyperLink"

$wb.Save()
